$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Gantt chart table updates (Report/ganttChart.xlsx) ---
# Column E (Days) is a shared formula (=D-C) and recalculates on its own
# once the Start/End dates change, so only the date cells are written.

# Row 9  "Functionality to flag patients": End slips a week (7 -> 14 days)
$ws.Range("D9").Value = 45372

# Row 10 "Generate reports/graphs": End slips a week (7 -> 14 days)
$ws.Range("D10").Value = 45372

# Row 11 "Testing": whole task shifts a week later (Start & End), stays 7 days
$ws.Range("C11").Value = 45372
$ws.Range("D11").Value = 45379

# --- Update the selected cell left behind by the editing session ---
$ws.Range("D17").Select()
